# Weekly price-sheet update: insert a new week's record for
# "Vega Modelo de Temuco - Ciboulette" at row 429, pushing the existing
# rows 429:449 down to 430:450 (dimension grows from R449 to R450).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 429; this shifts all rows
# from 429 down (old 429 -> 430, old 430 -> 431, ..., old 449 -> 450)
# and carries the date-column (D) number format down with it.
$ws.Rows.Item(429).Insert()

# Populate the newly inserted row 429 with this week's record.
$ws.Range("A429").Value = 10
$ws.Range("B429").Value = "Vega Modelo de Temuco"
$ws.Range("C429").Value = "La Araucanía"
$ws.Range("D429").Value = 45267
$ws.Range("E429").Value = 9
$ws.Range("F429").Value = 100112039
$ws.Range("G429").Value = "Ciboulette"
$ws.Range("H429").Value = "Sin especificar"
$ws.Range("I429").Value = "Primera"
$ws.Range("J429").Value = 40
$ws.Range("K429").Value = 7000
$ws.Range("L429").Value = 7000
$ws.Range("M429").Value = 7000
$ws.Range("N429").Value = "`$/docena de atados"
$ws.Range("O429").Value = "Provincia de Cautín"
$ws.Range("P429").Value = 2333
$ws.Range("Q429").Value = 3
$ws.Range("R429").Value = "Hortaliza"
